# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (94) currently carries the "latest day" number format
# (date-only, "YYYY-MM-DD"). Since it is no longer the latest day, it reverts
# to the regular datetime format used by every other historical row.
$ws.Range("A94").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row of data for the next day (day value 45682, i.e. 2025-01-25).
$newRow = 95
$ws.Cells.Item($newRow, 1).Value = 45682
$ws.Cells.Item($newRow, 2).Value = 227
$ws.Cells.Item($newRow, 3).Value = 223
$ws.Cells.Item($newRow, 4).Value = 221

# The new last row gets the "latest day" date-only number format that row 94
# used to have.
$ws.Range("A95").NumberFormat = "YYYY-MM-DD"
